# Renames the "_old" / "_new" header suffixes to the respective input-file
# format-version suffixes ("_FV2210" / "_FV2304"), wraps the sheet's data
# range in an Excel Table, and freezes the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row cells -----------------------------------
# Columns A-J carried the "_old" suffix, columns L-U carried "_new".
# Column K ("diff") is left untouched.
$oldHeaders = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210"
)

$newHeaders = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

for ($i = 0; $i -lt $oldHeaders.Count; $i++) {
    # Columns A..J -> indices 1..10
    $ws.Cells.Item(1, $i + 1).Value = $oldHeaders[$i]
}

for ($i = 0; $i -lt $newHeaders.Count; $i++) {
    # Columns L..U -> indices 12..21
    $ws.Cells.Item(1, $i + 12).Value = $newHeaders[$i]
}

# --- 2. Turn the used range into an Excel Table -------------------------
$dataRange = $ws.Range("A1:U57")
$lo = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$lo.Name = "Table1"

# --- 3. Freeze the header row -------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
